$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New resistance-calculation block in O4:X9 (shared strings must be created
# in this exact order to match the target uniqueCount/order).
# ---------------------------------------------------------------------------

# Row 4 - header banner (merged O4:X4)
$ws.Range("O4:X4").Merge()
$ws.Range("O4").Value = "The analogue readings are mappings of actual voltage into 1024 units"

# Row 5
$ws.Range("O5").Value = "Therefore:"
$ws.Range("P5:Q5").Merge()
$ws.Range("P5").Value = "5Vin/1024 = 4.9mV per unit"

# Row 6
$ws.Range("O6").Value = "meaning:"
$ws.Range("P6:Q6").Merge()
$ws.Range("P6").Value = "for example, lets take readign 40"
$ws.Range("R6:S6").Merge()
$ws.Range("R6").Value = "40*0.0049 = 0.196Vreading"

# Row 7
$ws.Range("O7").Value = "while:"
$ws.Range("P7:S7").Merge()
$ws.Range("P7").Value = "I = Vr1/R1 = (5-0.196)V/220ohm = 4.804Vr1/220ohm= 0.022A = 22mA"
$ws.Range("T7:V7").Merge()
$ws.Range("T7").Value = "using R1 to find current in circuit"

# Row 8
$ws.Range("O8").Value = "thus:"
$ws.Range("P8:R8").Merge()
$ws.Range("P8").Value = "R2 = Vreading/I = 0.196V/0.022 = 8.91ohm"

# Row 9 - footer banner (merged O9:X9)
$ws.Range("O9:X9").Merge()
$ws.Range("O9").Value = "*these calculations will be implemented within the code*"

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Full grid of thin box borders for the whole O5:X8 block (drawn first so the
# later, more specific borders on row 8 can override it).
$ws.Range("O5:X8").Borders.LineStyle = 1

# Header banner formatting: bold + underlined, centred both ways, box border.
$ws.Range("O4:X4").Font.Bold = $true
$ws.Range("O4:X4").Font.Underline = $true
$ws.Range("O4:X4").HorizontalAlignment = -4108
$ws.Range("O4:X4").VerticalAlignment = -4108
$ws.Range("O4:X4").Borders.LineStyle = 1
$ws.Range("O4:R4").Font.Size = 11
$ws.Range("S4:X4").Font.Size = 10

# Label column + value cells for rows 5-7: centred horizontally.
$ws.Range("O5:O8").HorizontalAlignment = -4108
$ws.Range("P5:Q5").HorizontalAlignment = -4108
$ws.Range("P6:Q6").HorizontalAlignment = -4108
$ws.Range("R6:S6").HorizontalAlignment = -4108
$ws.Range("P7:S7").HorizontalAlignment = -4108
$ws.Range("T7:V7").HorizontalAlignment = -4108

# Row 8 result band: continuous top/bottom border from P to V, with the box
# border only closed on the left (at P) and right (at V) edges.
$ws.Range("P8:V8").Borders.LineStyle = 0
$ws.Range("P8:V8").Borders(xlEdgeTop).LineStyle = 1
$ws.Range("P8:V8").Borders(xlEdgeBottom).LineStyle = 1
$ws.Range("P8").Borders(xlEdgeLeft).LineStyle = 1
$ws.Range("V8").Borders(xlEdgeRight).LineStyle = 1
$ws.Range("P8:R8").HorizontalAlignment = -4108
$ws.Range("W8:X8").Borders.LineStyle = 1

# Footer banner formatting: bold, centred horizontally, box border.
$ws.Range("O9:X9").Font.Bold = $true
$ws.Range("O9:X9").Font.Size = 11
$ws.Range("O9:X9").HorizontalAlignment = -4108
$ws.Range("O9:X9").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths for the new block.
# ---------------------------------------------------------------------------
$ws.Range("P1:R1").EntireColumn.ColumnWidth = 15
$ws.Range("S1:T1").EntireColumn.ColumnWidth = 13

# ---------------------------------------------------------------------------
# View state: scroll to the new block and select the P8:R8 result cell.
# ---------------------------------------------------------------------------
$ws.Range("P8:R8").Select()
$excel.ActiveWindow.ScrollColumn = 4
